$wb = $excel.ActiveWorkbook

# --- 1. Rename "Requested quantity" headers on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after "Monthly Trend" (at the end) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# --- 3. Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# reuse the same header style (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# reuse the same date-formatted style used for column A elsewhere
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A28").PasteSpecial(-4122)

$newSheet.Cells.Item(2,1).Value = 44990.99999999999
$newSheet.Cells.Item(2,2).Value = 24
$newSheet.Cells.Item(2,3).Value = -19.7158664843233
$newSheet.Cells.Item(2,4).Value = 67.16542472871372
$newSheet.Cells.Item(3,1).Value = 44997.99999999999
$newSheet.Cells.Item(3,2).Value = 25
$newSheet.Cells.Item(3,3).Value = -18.99572228883279
$newSheet.Cells.Item(3,4).Value = 67.13792218266335
$newSheet.Cells.Item(4,1).Value = 45004.99999999999
$newSheet.Cells.Item(4,2).Value = 25
$newSheet.Cells.Item(4,3).Value = -17.76603017529434
$newSheet.Cells.Item(4,4).Value = 68.09339095020984
$newSheet.Cells.Item(5,1).Value = 45011.99999999999
$newSheet.Cells.Item(5,2).Value = 25
$newSheet.Cells.Item(5,3).Value = -18.24042704514262
$newSheet.Cells.Item(5,4).Value = 67.17618880475422
$newSheet.Cells.Item(6,1).Value = 45039.99999999999
$newSheet.Cells.Item(6,2).Value = 26
$newSheet.Cells.Item(6,3).Value = -19.38011486304048
$newSheet.Cells.Item(6,4).Value = 69.75232174194585
$newSheet.Cells.Item(7,1).Value = 45088.99999999999
$newSheet.Cells.Item(7,2).Value = 27
$newSheet.Cells.Item(7,3).Value = -16.43882613197266
$newSheet.Cells.Item(7,4).Value = 70.78863551678161
$newSheet.Cells.Item(8,1).Value = 45179.99999999999
$newSheet.Cells.Item(8,2).Value = 30
$newSheet.Cells.Item(8,3).Value = -15.62413630091297
$newSheet.Cells.Item(8,4).Value = 73.04120487178982
$newSheet.Cells.Item(9,1).Value = 45186.99999999999
$newSheet.Cells.Item(9,2).Value = 30
$newSheet.Cells.Item(9,3).Value = -12.51698442554778
$newSheet.Cells.Item(9,4).Value = 74.18645627003768
$newSheet.Cells.Item(10,1).Value = 45193.99999999999
$newSheet.Cells.Item(10,2).Value = 31
$newSheet.Cells.Item(10,3).Value = -14.26574581603047
$newSheet.Cells.Item(10,4).Value = 77.32793492744558
$newSheet.Cells.Item(11,1).Value = 45200.99999999999
$newSheet.Cells.Item(11,2).Value = 31
$newSheet.Cells.Item(11,3).Value = -12.91545283230738
$newSheet.Cells.Item(11,4).Value = 74.00493938136241
$newSheet.Cells.Item(12,1).Value = 45214.99999999999
$newSheet.Cells.Item(12,2).Value = 31
$newSheet.Cells.Item(12,3).Value = -12.70755795373577
$newSheet.Cells.Item(12,4).Value = 71.82949615960599
$newSheet.Cells.Item(13,1).Value = 45221.99999999999
$newSheet.Cells.Item(13,2).Value = 32
$newSheet.Cells.Item(13,3).Value = -11.34820661982495
$newSheet.Cells.Item(13,4).Value = 74.15787693841337
$newSheet.Cells.Item(14,1).Value = 45228.99999999999
$newSheet.Cells.Item(14,2).Value = 32
$newSheet.Cells.Item(14,3).Value = -10.08592353450919
$newSheet.Cells.Item(14,4).Value = 76.39589642272736
$newSheet.Cells.Item(15,1).Value = 45235.99999999999
$newSheet.Cells.Item(15,2).Value = 32
$newSheet.Cells.Item(15,3).Value = -13.59862143221971
$newSheet.Cells.Item(15,4).Value = 75.42247980902873
$newSheet.Cells.Item(16,1).Value = 45242.99999999999
$newSheet.Cells.Item(16,2).Value = 32
$newSheet.Cells.Item(16,3).Value = -11.14477219577557
$newSheet.Cells.Item(16,4).Value = 74.76423138825974
$newSheet.Cells.Item(17,1).Value = 45249.99999999999
$newSheet.Cells.Item(17,2).Value = 32
$newSheet.Cells.Item(17,3).Value = -8.667329815069493
$newSheet.Cells.Item(17,4).Value = 77.14738994749501
$newSheet.Cells.Item(18,1).Value = 45270.99999999999
$newSheet.Cells.Item(18,2).Value = 33
$newSheet.Cells.Item(18,3).Value = -11.2669007853591
$newSheet.Cells.Item(18,4).Value = 75.9689163751552
$newSheet.Cells.Item(19,1).Value = 45613.99999999999
$newSheet.Cells.Item(19,2).Value = 44
$newSheet.Cells.Item(19,3).Value = -0.4048742598760813
$newSheet.Cells.Item(19,4).Value = 83.96287926987776
$newSheet.Cells.Item(20,1).Value = 45641.99999999999
$newSheet.Cells.Item(20,2).Value = 45
$newSheet.Cells.Item(20,3).Value = 3.01272273653262
$newSheet.Cells.Item(20,4).Value = 87.47160815754653
$newSheet.Cells.Item(21,1).Value = 45648.99999999999
$newSheet.Cells.Item(21,2).Value = 45
$newSheet.Cells.Item(21,3).Value = -0.7139054043507204
$newSheet.Cells.Item(21,4).Value = 87.40300316166268
$newSheet.Cells.Item(22,1).Value = 45655.99999999999
$newSheet.Cells.Item(22,2).Value = 45
$newSheet.Cells.Item(22,3).Value = 2.452761518380349
$newSheet.Cells.Item(22,4).Value = 87.97093590506938
$newSheet.Cells.Item(23,1).Value = 45662.99999999999
$newSheet.Cells.Item(23,2).Value = 45
$newSheet.Cells.Item(23,3).Value = 2.174642053432746
$newSheet.Cells.Item(23,4).Value = 89.01684690420548
$newSheet.Cells.Item(24,1).Value = 45669.99999999999
$newSheet.Cells.Item(24,2).Value = 46
$newSheet.Cells.Item(24,3).Value = 2.083206700876225
$newSheet.Cells.Item(24,4).Value = 86.2665212536115
$newSheet.Cells.Item(25,1).Value = 45676.99999999999
$newSheet.Cells.Item(25,2).Value = 46
$newSheet.Cells.Item(25,3).Value = 2.901309935539547
$newSheet.Cells.Item(25,4).Value = 89.02139390584776
$newSheet.Cells.Item(26,1).Value = 45683.99999999999
$newSheet.Cells.Item(26,2).Value = 46
$newSheet.Cells.Item(26,3).Value = 6.989686417798728
$newSheet.Cells.Item(26,4).Value = 87.36375072997629
$newSheet.Cells.Item(27,1).Value = 45690.99999999999
$newSheet.Cells.Item(27,2).Value = 46
$newSheet.Cells.Item(27,3).Value = 2.927651652246611
$newSheet.Cells.Item(27,4).Value = 87.56647361621589
$newSheet.Cells.Item(28,1).Value = 45697.99999999999
$newSheet.Cells.Item(28,2).Value = 47
$newSheet.Cells.Item(28,3).Value = -0.02061590056212335
$newSheet.Cells.Item(28,4).Value = 88.76452284197958
